$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) counts replacing the old Strike# counts in column G
$kValues = @(7,8,3,6,5,6,7,7,11,8,4,6,8,8,8,4,7,6,5,4,9,7,6,11,5,10,7,0,7,6,7,5,5,6,1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
